$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 252 (previously did not exist) - full record, shifted from old row 251
$ws.Range("A252").Value = 10
$ws.Range("B252").Value = 'Vega Modelo de Temuco'
$ws.Range("C252").Value = 'La Araucanía'
$ws.Range("D252").Value = 44487
$ws.Range("E252").Value = 9
$ws.Range("F252").Value = 100112039
$ws.Range("G252").Value = 'Ciboulette'
$ws.Range("H252").Value = 'Sin especificar'
$ws.Range("I252").Value = 'Primera'
$ws.Range("J252").Value = 75
$ws.Range("K252").Value = 2000
$ws.Range("L252").Value = 2000
$ws.Range("M252").Value = 2000
$ws.Range("N252").Value = '$/docena de atados'
$ws.Range("O252").Value = 'Región Metropolitana'
$ws.Range("P252").Value = 667
$ws.Range("Q252").Value = 3
$ws.Range("R252").Value = 'Hortaliza'
$ws.Range("D252").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update rows 158-251: row 158 becomes a new record; rows 159-251 shift down from the previous row
$ws.Range("D158").Value = 44719
$ws.Range("J158").Value = 20
$ws.Range("K158").Value = 6000
$ws.Range("L158").Value = 6000
$ws.Range("M158").Value = 6000
$ws.Range("O158").Value = 'Provincia de Cautín'
$ws.Range("P158").Value = 2000
$ws.Range("D159").Value = 44637
$ws.Range("J159").Value = 40
$ws.Range("K159").Value = 5000
$ws.Range("L159").Value = 5000
$ws.Range("M159").Value = 5000
$ws.Range("O159").Value = 'Provincia de Cautín'
$ws.Range("P159").Value = 1667
$ws.Range("D160").Value = 44461
$ws.Range("J160").Value = 40
$ws.Range("K160").Value = 3000
$ws.Range("L160").Value = 6000
$ws.Range("M160").Value = 3750
$ws.Range("O160").Value = 'Provincia de Cautín'
$ws.Range("P160").Value = 1250
$ws.Range("D161").Value = 44600
$ws.Range("J161").Value = 110
$ws.Range("K161").Value = 5000
$ws.Range("L161").Value = 5000
$ws.Range("M161").Value = 5000
$ws.Range("O161").Value = 'Provincia de Cautín'
$ws.Range("P161").Value = 1667
$ws.Range("D162").Value = 44412
$ws.Range("J162").Value = 30
$ws.Range("K162").Value = 5000
$ws.Range("L162").Value = 5000
$ws.Range("M162").Value = 5000
$ws.Range("O162").Value = 'Región Metropolitana'
$ws.Range("P162").Value = 1667
$ws.Range("D163").Value = 44302
$ws.Range("J163").Value = 20
$ws.Range("K163").Value = 7000
$ws.Range("L163").Value = 7000
$ws.Range("M163").Value = 7000
$ws.Range("O163").Value = 'Provincia de Cautín'
$ws.Range("P163").Value = 2333
$ws.Range("D164").Value = 44477
$ws.Range("J164").Value = 20
$ws.Range("K164").Value = 7000
$ws.Range("L164").Value = 7000
$ws.Range("M164").Value = 7000
$ws.Range("O164").Value = 'Provincia de Cautín'
$ws.Range("P164").Value = 2333
$ws.Range("D165").Value = 44523
$ws.Range("J165").Value = 30
$ws.Range("K165").Value = 5000
$ws.Range("L165").Value = 5000
$ws.Range("M165").Value = 5000
$ws.Range("O165").Value = 'Provincia de Cautín'
$ws.Range("P165").Value = 1667
$ws.Range("D166").Value = 44396
$ws.Range("J166").Value = 30
$ws.Range("K166").Value = 7000
$ws.Range("L166").Value = 7000
$ws.Range("M166").Value = 7000
$ws.Range("O166").Value = 'Provincia de Cautín'
$ws.Range("P166").Value = 2333
$ws.Range("D167").Value = 44364
$ws.Range("J167").Value = 65
$ws.Range("K167").Value = 2500
$ws.Range("L167").Value = 2500
$ws.Range("M167").Value = 2500
$ws.Range("O167").Value = 'Región Metropolitana'
$ws.Range("P167").Value = 833
$ws.Range("D168").Value = 44434
$ws.Range("J168").Value = 45
$ws.Range("K168").Value = 8000
$ws.Range("L168").Value = 8000
$ws.Range("M168").Value = 8000
$ws.Range("O168").Value = 'Provincia de Cautín'
$ws.Range("P168").Value = 2667
$ws.Range("D169").Value = 44400
$ws.Range("J169").Value = 10
$ws.Range("K169").Value = 10000
$ws.Range("L169").Value = 10000
$ws.Range("M169").Value = 10000
$ws.Range("O169").Value = 'Provincia de Cautín'
$ws.Range("P169").Value = 3333
$ws.Range("D170").Value = 44413
$ws.Range("J170").Value = 20
$ws.Range("K170").Value = 10000
$ws.Range("L170").Value = 10000
$ws.Range("M170").Value = 10000
$ws.Range("O170").Value = 'Provincia de Cautín'
$ws.Range("P170").Value = 3333
$ws.Range("D171").Value = 44649
$ws.Range("J171").Value = 40
$ws.Range("K171").Value = 5000
$ws.Range("L171").Value = 5000
$ws.Range("M171").Value = 5000
$ws.Range("O171").Value = 'Provincia de Cautín'
$ws.Range("P171").Value = 1667
$ws.Range("D172").Value = 44406
$ws.Range("J172").Value = 120
$ws.Range("K172").Value = 9000
$ws.Range("L172").Value = 10000
$ws.Range("M172").Value = 9542
$ws.Range("O172").Value = 'Provincia de Cautín'
$ws.Range("P172").Value = 3181
$ws.Range("D173").Value = 44432
$ws.Range("J173").Value = 15
$ws.Range("K173").Value = 7000
$ws.Range("L173").Value = 7000
$ws.Range("M173").Value = 7000
$ws.Range("O173").Value = 'Provincia de Cautín'
$ws.Range("P173").Value = 2333
$ws.Range("D174").Value = 44315
$ws.Range("J174").Value = 40
$ws.Range("K174").Value = 6000
$ws.Range("L174").Value = 6000
$ws.Range("M174").Value = 6000
$ws.Range("O174").Value = 'Provincia de Cautín'
$ws.Range("P174").Value = 2000
$ws.Range("D175").Value = 44575
$ws.Range("J175").Value = 65
$ws.Range("K175").Value = 5000
$ws.Range("L175").Value = 5000
$ws.Range("M175").Value = 5000
$ws.Range("O175").Value = 'Provincia de Cautín'
$ws.Range("P175").Value = 1667
$ws.Range("D176").Value = 44257
$ws.Range("J176").Value = 30
$ws.Range("K176").Value = 5000
$ws.Range("L176").Value = 5000
$ws.Range("M176").Value = 5000
$ws.Range("O176").Value = 'Provincia de Cautín'
$ws.Range("P176").Value = 1667
$ws.Range("D177").Value = 44428
$ws.Range("J177").Value = 10
$ws.Range("K177").Value = 8000
$ws.Range("L177").Value = 8000
$ws.Range("M177").Value = 8000
$ws.Range("O177").Value = 'Provincia de Cautín'
$ws.Range("P177").Value = 2667
$ws.Range("D178").Value = 44596
$ws.Range("J178").Value = 20
$ws.Range("K178").Value = 5000
$ws.Range("L178").Value = 5000
$ws.Range("M178").Value = 5000
$ws.Range("O178").Value = 'Provincia de Cautín'
$ws.Range("P178").Value = 1667
$ws.Range("D179").Value = 44536
$ws.Range("J179").Value = 65
$ws.Range("K179").Value = 5000
$ws.Range("L179").Value = 5000
$ws.Range("M179").Value = 5000
$ws.Range("O179").Value = 'Provincia de Cautín'
$ws.Range("P179").Value = 1667
$ws.Range("D180").Value = 44350
$ws.Range("J180").Value = 55
$ws.Range("K180").Value = 6000
$ws.Range("L180").Value = 6000
$ws.Range("M180").Value = 6000
$ws.Range("O180").Value = 'Provincia de Cautín'
$ws.Range("P180").Value = 2000
$ws.Range("D181").Value = 44630
$ws.Range("J181").Value = 30
$ws.Range("K181").Value = 5000
$ws.Range("L181").Value = 5000
$ws.Range("M181").Value = 5000
$ws.Range("O181").Value = 'Provincia de Cautín'
$ws.Range("P181").Value = 1667
$ws.Range("D182").Value = 44467
$ws.Range("J182").Value = 20
$ws.Range("K182").Value = 7000
$ws.Range("L182").Value = 7000
$ws.Range("M182").Value = 7000
$ws.Range("O182").Value = 'Provincia de Cautín'
$ws.Range("P182").Value = 2333
$ws.Range("D183").Value = 44448
$ws.Range("J183").Value = 65
$ws.Range("K183").Value = 8000
$ws.Range("L183").Value = 8000
$ws.Range("M183").Value = 8000
$ws.Range("O183").Value = 'Provincia de Cautín'
$ws.Range("P183").Value = 2667
$ws.Range("D184").Value = 44420
$ws.Range("J184").Value = 120
$ws.Range("K184").Value = 10000
$ws.Range("L184").Value = 10000
$ws.Range("M184").Value = 10000
$ws.Range("O184").Value = 'Provincia de Cautín'
$ws.Range("P184").Value = 3333
$ws.Range("D185").Value = 44308
$ws.Range("J185").Value = 65
$ws.Range("K185").Value = 6000
$ws.Range("L185").Value = 7000
$ws.Range("M185").Value = 6538
$ws.Range("O185").Value = 'Provincia de Cautín'
$ws.Range("P185").Value = 2179
$ws.Range("D186").Value = 44435
$ws.Range("J186").Value = 180
$ws.Range("K186").Value = 7000
$ws.Range("L186").Value = 8000
$ws.Range("M186").Value = 7917
$ws.Range("O186").Value = 'Provincia de Cautín'
$ws.Range("P186").Value = 2639
$ws.Range("D187").Value = 44208
$ws.Range("J187").Value = 65
$ws.Range("K187").Value = 4000
$ws.Range("L187").Value = 4000
$ws.Range("M187").Value = 4000
$ws.Range("O187").Value = 'Provincia de Cautín'
$ws.Range("P187").Value = 1333
$ws.Range("D188").Value = 44224
$ws.Range("J188").Value = 100
$ws.Range("K188").Value = 5000
$ws.Range("L188").Value = 6000
$ws.Range("M188").Value = 5550
$ws.Range("O188").Value = 'Provincia de Cautín'
$ws.Range("P188").Value = 1850
$ws.Range("D189").Value = 44259
$ws.Range("J189").Value = 80
$ws.Range("K189").Value = 5000
$ws.Range("L189").Value = 5000
$ws.Range("M189").Value = 5000
$ws.Range("O189").Value = 'Provincia de Cautín'
$ws.Range("P189").Value = 1667
$ws.Range("D190").Value = 44449
$ws.Range("J190").Value = 95
$ws.Range("K190").Value = 8000
$ws.Range("L190").Value = 8000
$ws.Range("M190").Value = 8000
$ws.Range("O190").Value = 'Provincia de Cautín'
$ws.Range("P190").Value = 2667
$ws.Range("D191").Value = 44449
$ws.Range("J191").Value = 85
$ws.Range("K191").Value = 4000
$ws.Range("L191").Value = 4000
$ws.Range("M191").Value = 4000
$ws.Range("O191").Value = 'Región Metropolitana'
$ws.Range("P191").Value = 1333
$ws.Range("D192").Value = 44410
$ws.Range("J192").Value = 30
$ws.Range("K192").Value = 8000
$ws.Range("L192").Value = 8000
$ws.Range("M192").Value = 8000
$ws.Range("O192").Value = 'Provincia de Cautín'
$ws.Range("P192").Value = 2667
$ws.Range("D193").Value = 44508
$ws.Range("J193").Value = 40
$ws.Range("K193").Value = 5000
$ws.Range("L193").Value = 5000
$ws.Range("M193").Value = 5000
$ws.Range("O193").Value = 'Provincia de Cautín'
$ws.Range("P193").Value = 1667
$ws.Range("D194").Value = 44571
$ws.Range("J194").Value = 65
$ws.Range("K194").Value = 5000
$ws.Range("L194").Value = 5000
$ws.Range("M194").Value = 5000
$ws.Range("O194").Value = 'Provincia de Cautín'
$ws.Range("P194").Value = 1667
$ws.Range("D195").Value = 44512
$ws.Range("J195").Value = 20
$ws.Range("K195").Value = 5000
$ws.Range("L195").Value = 5000
$ws.Range("M195").Value = 5000
$ws.Range("O195").Value = 'Provincia de Cautín'
$ws.Range("P195").Value = 1667
$ws.Range("D196").Value = 44490
$ws.Range("J196").Value = 65
$ws.Range("K196").Value = 6000
$ws.Range("L196").Value = 6000
$ws.Range("M196").Value = 6000
$ws.Range("O196").Value = 'Provincia de Cautín'
$ws.Range("P196").Value = 2000
$ws.Range("D197").Value = 44418
$ws.Range("J197").Value = 65
$ws.Range("K197").Value = 10000
$ws.Range("L197").Value = 10000
$ws.Range("M197").Value = 10000
$ws.Range("O197").Value = 'Provincia de Cautín'
$ws.Range("P197").Value = 3333
$ws.Range("D198").Value = 44203
$ws.Range("J198").Value = 80
$ws.Range("K198").Value = 4000
$ws.Range("L198").Value = 4000
$ws.Range("M198").Value = 4000
$ws.Range("O198").Value = 'Provincia de Cautín'
$ws.Range("P198").Value = 1333
$ws.Range("D199").Value = 44679
$ws.Range("J199").Value = 70
$ws.Range("K199").Value = 6000
$ws.Range("L199").Value = 8000
$ws.Range("M199").Value = 7143
$ws.Range("O199").Value = 'Provincia de Cautín'
$ws.Range("P199").Value = 2381
$ws.Range("D200").Value = 44663
$ws.Range("J200").Value = 30
$ws.Range("K200").Value = 5000
$ws.Range("L200").Value = 5000
$ws.Range("M200").Value = 5000
$ws.Range("O200").Value = 'Provincia de Cautín'
$ws.Range("P200").Value = 1667
$ws.Range("D201").Value = 44174
$ws.Range("J201").Value = 30
$ws.Range("K201").Value = 4000
$ws.Range("L201").Value = 4000
$ws.Range("M201").Value = 4000
$ws.Range("O201").Value = 'Provincia de Cautín'
$ws.Range("P201").Value = 1333
$ws.Range("D202").Value = 44424
$ws.Range("J202").Value = 20
$ws.Range("K202").Value = 8000
$ws.Range("L202").Value = 8000
$ws.Range("M202").Value = 8000
$ws.Range("O202").Value = 'Provincia de Cautín'
$ws.Range("P202").Value = 2667
$ws.Range("D203").Value = 44623
$ws.Range("J203").Value = 40
$ws.Range("K203").Value = 5000
$ws.Range("L203").Value = 5000
$ws.Range("M203").Value = 5000
$ws.Range("O203").Value = 'Provincia de Cautín'
$ws.Range("P203").Value = 1667
$ws.Range("D204").Value = 44313
$ws.Range("J204").Value = 30
$ws.Range("K204").Value = 5000
$ws.Range("L204").Value = 6000
$ws.Range("M204").Value = 5333
$ws.Range("O204").Value = 'Provincia de Cautín'
$ws.Range("P204").Value = 1778
$ws.Range("D205").Value = 44664
$ws.Range("J205").Value = 30
$ws.Range("K205").Value = 5000
$ws.Range("L205").Value = 5000
$ws.Range("M205").Value = 5000
$ws.Range("O205").Value = 'Provincia de Cautín'
$ws.Range("P205").Value = 1667
$ws.Range("D206").Value = 44195
$ws.Range("J206").Value = 55
$ws.Range("K206").Value = 6000
$ws.Range("L206").Value = 6000
$ws.Range("M206").Value = 6000
$ws.Range("O206").Value = 'Provincia de Cautín'
$ws.Range("P206").Value = 2000
$ws.Range("D207").Value = 44274
$ws.Range("J207").Value = 20
$ws.Range("K207").Value = 5000
$ws.Range("L207").Value = 5000
$ws.Range("M207").Value = 5000
$ws.Range("O207").Value = 'Provincia de Cautín'
$ws.Range("P207").Value = 1667
$ws.Range("D208").Value = 44417
$ws.Range("J208").Value = 65
$ws.Range("K208").Value = 10000
$ws.Range("L208").Value = 10000
$ws.Range("M208").Value = 10000
$ws.Range("O208").Value = 'Provincia de Cautín'
$ws.Range("P208").Value = 3333
$ws.Range("D209").Value = 44578
$ws.Range("J209").Value = 50
$ws.Range("K209").Value = 5000
$ws.Range("L209").Value = 5000
$ws.Range("M209").Value = 5000
$ws.Range("O209").Value = 'Provincia de Cautín'
$ws.Range("P209").Value = 1667
$ws.Range("D210").Value = 44495
$ws.Range("J210").Value = 20
$ws.Range("K210").Value = 6000
$ws.Range("L210").Value = 6000
$ws.Range("M210").Value = 6000
$ws.Range("O210").Value = 'Provincia de Cautín'
$ws.Range("P210").Value = 2000
$ws.Range("D211").Value = 44648
$ws.Range("J211").Value = 30
$ws.Range("K211").Value = 5000
$ws.Range("L211").Value = 5000
$ws.Range("M211").Value = 5000
$ws.Range("O211").Value = 'Provincia de Cautín'
$ws.Range("P211").Value = 1667
$ws.Range("D212").Value = 44221
$ws.Range("J212").Value = 85
$ws.Range("K212").Value = 5000
$ws.Range("L212").Value = 6000
$ws.Range("M212").Value = 5588
$ws.Range("O212").Value = 'Provincia de Cautín'
$ws.Range("P212").Value = 1863
$ws.Range("D213").Value = 44580
$ws.Range("J213").Value = 40
$ws.Range("K213").Value = 5000
$ws.Range("L213").Value = 5000
$ws.Range("M213").Value = 5000
$ws.Range("O213").Value = 'Provincia de Cautín'
$ws.Range("P213").Value = 1667
$ws.Range("D214").Value = 44371
$ws.Range("J214").Value = 40
$ws.Range("K214").Value = 5000
$ws.Range("L214").Value = 5000
$ws.Range("M214").Value = 5000
$ws.Range("O214").Value = 'Provincia de Cautín'
$ws.Range("P214").Value = 1667
$ws.Range("D215").Value = 44708
$ws.Range("J215").Value = 20
$ws.Range("K215").Value = 6000
$ws.Range("L215").Value = 6000
$ws.Range("M215").Value = 6000
$ws.Range("O215").Value = 'Provincia de Cautín'
$ws.Range("P215").Value = 2000
$ws.Range("D216").Value = 44454
$ws.Range("J216").Value = 20
$ws.Range("K216").Value = 8000
$ws.Range("L216").Value = 8000
$ws.Range("M216").Value = 8000
$ws.Range("O216").Value = 'Provincia de Cautín'
$ws.Range("P216").Value = 2667
$ws.Range("D217").Value = 44160
$ws.Range("J217").Value = 20
$ws.Range("K217").Value = 5000
$ws.Range("L217").Value = 5000
$ws.Range("M217").Value = 5000
$ws.Range("O217").Value = 'Provincia de Cautín'
$ws.Range("P217").Value = 1667
$ws.Range("D218").Value = 44565
$ws.Range("J218").Value = 30
$ws.Range("K218").Value = 5000
$ws.Range("L218").Value = 5000
$ws.Range("M218").Value = 5000
$ws.Range("O218").Value = 'Provincia de Cautín'
$ws.Range("P218").Value = 1667
$ws.Range("D219").Value = 44603
$ws.Range("J219").Value = 65
$ws.Range("K219").Value = 5000
$ws.Range("L219").Value = 5000
$ws.Range("M219").Value = 5000
$ws.Range("O219").Value = 'Provincia de Cautín'
$ws.Range("P219").Value = 1667
$ws.Range("D220").Value = 44263
$ws.Range("J220").Value = 55
$ws.Range("K220").Value = 5000
$ws.Range("L220").Value = 5000
$ws.Range("M220").Value = 5000
$ws.Range("O220").Value = 'Provincia de Cautín'
$ws.Range("P220").Value = 1667
$ws.Range("D221").Value = 44187
$ws.Range("J221").Value = 30
$ws.Range("K221").Value = 4000
$ws.Range("L221").Value = 4000
$ws.Range("M221").Value = 4000
$ws.Range("O221").Value = 'Provincia de Cautín'
$ws.Range("P221").Value = 1333
$ws.Range("D222").Value = 44609
$ws.Range("J222").Value = 40
$ws.Range("K222").Value = 5000
$ws.Range("L222").Value = 5000
$ws.Range("M222").Value = 5000
$ws.Range("O222").Value = 'Provincia de Cautín'
$ws.Range("P222").Value = 1667
$ws.Range("D223").Value = 44529
$ws.Range("J223").Value = 110
$ws.Range("K223").Value = 5000
$ws.Range("L223").Value = 5000
$ws.Range("M223").Value = 5000
$ws.Range("O223").Value = 'Provincia de Cautín'
$ws.Range("P223").Value = 1667
$ws.Range("D224").Value = 44321
$ws.Range("J224").Value = 45
$ws.Range("K224").Value = 6000
$ws.Range("L224").Value = 6000
$ws.Range("M224").Value = 6000
$ws.Range("O224").Value = 'Provincia de Cautín'
$ws.Range("P224").Value = 2000
$ws.Range("D225").Value = 44277
$ws.Range("J225").Value = 65
$ws.Range("K225").Value = 4000
$ws.Range("L225").Value = 4000
$ws.Range("M225").Value = 4000
$ws.Range("O225").Value = 'Provincia de Cautín'
$ws.Range("P225").Value = 1333
$ws.Range("D226").Value = 44166
$ws.Range("J226").Value = 65
$ws.Range("K226").Value = 5000
$ws.Range("L226").Value = 5000
$ws.Range("M226").Value = 5000
$ws.Range("O226").Value = 'Provincia de Cautín'
$ws.Range("P226").Value = 1667
$ws.Range("D227").Value = 44354
$ws.Range("J227").Value = 30
$ws.Range("K227").Value = 5000
$ws.Range("L227").Value = 5000
$ws.Range("M227").Value = 5000
$ws.Range("O227").Value = 'Provincia de Cautín'
$ws.Range("P227").Value = 1667
$ws.Range("D228").Value = 44245
$ws.Range("J228").Value = 65
$ws.Range("K228").Value = 5000
$ws.Range("L228").Value = 5000
$ws.Range("M228").Value = 5000
$ws.Range("O228").Value = 'Provincia de Cautín'
$ws.Range("P228").Value = 1667
$ws.Range("D229").Value = 44168
$ws.Range("J229").Value = 85
$ws.Range("K229").Value = 5000
$ws.Range("L229").Value = 5000
$ws.Range("M229").Value = 5000
$ws.Range("O229").Value = 'Provincia de Cautín'
$ws.Range("P229").Value = 1667
$ws.Range("D230").Value = 44638
$ws.Range("J230").Value = 30
$ws.Range("K230").Value = 5000
$ws.Range("L230").Value = 5000
$ws.Range("M230").Value = 5000
$ws.Range("O230").Value = 'Provincia de Cautín'
$ws.Range("P230").Value = 1667
$ws.Range("D231").Value = 44249
$ws.Range("J231").Value = 110
$ws.Range("K231").Value = 5000
$ws.Range("L231").Value = 5000
$ws.Range("M231").Value = 5000
$ws.Range("O231").Value = 'Provincia de Cautín'
$ws.Range("P231").Value = 1667
$ws.Range("D232").Value = 44431
$ws.Range("J232").Value = 55
$ws.Range("K232").Value = 8000
$ws.Range("L232").Value = 8000
$ws.Range("M232").Value = 8000
$ws.Range("O232").Value = 'Provincia de Cautín'
$ws.Range("P232").Value = 2667
$ws.Range("D233").Value = 44522
$ws.Range("J233").Value = 30
$ws.Range("K233").Value = 5000
$ws.Range("L233").Value = 5000
$ws.Range("M233").Value = 5000
$ws.Range("O233").Value = 'Provincia de Cautín'
$ws.Range("P233").Value = 1667
$ws.Range("D234").Value = 44714
$ws.Range("J234").Value = 65
$ws.Range("K234").Value = 6000
$ws.Range("L234").Value = 6000
$ws.Range("M234").Value = 6000
$ws.Range("O234").Value = 'Provincia de Cautín'
$ws.Range("P234").Value = 2000
$ws.Range("D235").Value = 44188
$ws.Range("J235").Value = 40
$ws.Range("K235").Value = 4000
$ws.Range("L235").Value = 5000
$ws.Range("M235").Value = 4500
$ws.Range("O235").Value = 'Provincia de Cautín'
$ws.Range("P235").Value = 1500
$ws.Range("D236").Value = 44659
$ws.Range("J236").Value = 65
$ws.Range("K236").Value = 6000
$ws.Range("L236").Value = 6000
$ws.Range("M236").Value = 6000
$ws.Range("O236").Value = 'Provincia de Cautín'
$ws.Range("P236").Value = 2000
$ws.Range("D237").Value = 44494
$ws.Range("J237").Value = 20
$ws.Range("K237").Value = 6000
$ws.Range("L237").Value = 6000
$ws.Range("M237").Value = 6000
$ws.Range("O237").Value = 'Provincia de Cautín'
$ws.Range("P237").Value = 2000
$ws.Range("D238").Value = 44494
$ws.Range("J238").Value = 30
$ws.Range("K238").Value = 2000
$ws.Range("L238").Value = 2000
$ws.Range("M238").Value = 2000
$ws.Range("O238").Value = 'Región Metropolitana'
$ws.Range("P238").Value = 667
$ws.Range("D239").Value = 44323
$ws.Range("J239").Value = 50
$ws.Range("K239").Value = 5000
$ws.Range("L239").Value = 5000
$ws.Range("M239").Value = 5000
$ws.Range("O239").Value = 'Provincia de Cautín'
$ws.Range("P239").Value = 1667
$ws.Range("D240").Value = 44526
$ws.Range("J240").Value = 20
$ws.Range("K240").Value = 5000
$ws.Range("L240").Value = 5000
$ws.Range("M240").Value = 5000
$ws.Range("O240").Value = 'Provincia de Cautín'
$ws.Range("P240").Value = 1667
$ws.Range("D241").Value = 44363
$ws.Range("J241").Value = 80
$ws.Range("K241").Value = 2500
$ws.Range("L241").Value = 2500
$ws.Range("M241").Value = 2500
$ws.Range("O241").Value = 'Región Metropolitana'
$ws.Range("P241").Value = 833
$ws.Range("D242").Value = 44704
$ws.Range("J242").Value = 30
$ws.Range("K242").Value = 6000
$ws.Range("L242").Value = 6000
$ws.Range("M242").Value = 6000
$ws.Range("O242").Value = 'Provincia de Cautín'
$ws.Range("P242").Value = 2000
$ws.Range("D243").Value = 44620
$ws.Range("J243").Value = 30
$ws.Range("K243").Value = 5000
$ws.Range("L243").Value = 5000
$ws.Range("M243").Value = 5000
$ws.Range("O243").Value = 'Provincia de Cautín'
$ws.Range("P243").Value = 1667
$ws.Range("D244").Value = 44586
$ws.Range("J244").Value = 65
$ws.Range("K244").Value = 5000
$ws.Range("L244").Value = 5000
$ws.Range("M244").Value = 5000
$ws.Range("O244").Value = 'Provincia de Cautín'
$ws.Range("P244").Value = 1667
$ws.Range("D245").Value = 44601
$ws.Range("J245").Value = 65
$ws.Range("K245").Value = 5000
$ws.Range("L245").Value = 5000
$ws.Range("M245").Value = 5000
$ws.Range("O245").Value = 'Provincia de Cautín'
$ws.Range("P245").Value = 1667
$ws.Range("D246").Value = 44544
$ws.Range("J246").Value = 25
$ws.Range("K246").Value = 7000
$ws.Range("L246").Value = 7000
$ws.Range("M246").Value = 7000
$ws.Range("O246").Value = 'Provincia de Cautín'
$ws.Range("P246").Value = 2333
$ws.Range("D247").Value = 44617
$ws.Range("J247").Value = 40
$ws.Range("K247").Value = 5000
$ws.Range("L247").Value = 5000
$ws.Range("M247").Value = 5000
$ws.Range("O247").Value = 'Provincia de Cautín'
$ws.Range("P247").Value = 1667
$ws.Range("D248").Value = 44567
$ws.Range("J248").Value = 40
$ws.Range("K248").Value = 5000
$ws.Range("L248").Value = 5000
$ws.Range("M248").Value = 5000
$ws.Range("O248").Value = 'Provincia de Cautín'
$ws.Range("P248").Value = 1667
$ws.Range("D249").Value = 44169
$ws.Range("J249").Value = 75
$ws.Range("K249").Value = 5000
$ws.Range("L249").Value = 5000
$ws.Range("M249").Value = 5000
$ws.Range("O249").Value = 'Provincia de Cautín'
$ws.Range("P249").Value = 1667
$ws.Range("D250").Value = 44474
$ws.Range("J250").Value = 30
$ws.Range("K250").Value = 5000
$ws.Range("L250").Value = 5000
$ws.Range("M250").Value = 5000
$ws.Range("O250").Value = 'Región Metropolitana'
$ws.Range("P250").Value = 1667
$ws.Range("D251").Value = 44487
$ws.Range("J251").Value = 105
$ws.Range("K251").Value = 5000
$ws.Range("L251").Value = 6000
$ws.Range("M251").Value = 5619
$ws.Range("O251").Value = 'Provincia de Cautín'
$ws.Range("P251").Value = 1873
